$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Tenure in years" column header (C1) is being removed, so the headers
# that followed it (D1:M1) all shift one column to the left (C1:L1), and the
# vacated M1 becomes empty. The column width definitions in <cols> are left
# completely untouched (this was a cell-content move, not a column delete),
# so we move the header cells individually rather than using Delete/Insert
# on whole columns (which also renumbers <cols>).
#
# Range.Cut()/Copy() with overlapping source/destination ranges behaves
# unreliably here, so the D1:M1 block is staged through a scratch row first
# (values, then formats, copied separately since PasteSpecial only honors
# one aspect at a time) and then moved from the scratch row into C1:L1.

$src     = $ws.Range("D1:M1")
$scratch = $ws.Range("D100:M100")
$dst     = $ws.Range("C1:L1")

$src.Copy()
$scratch.PasteSpecial(-4163)   # xlPasteValues
$src.Copy()
$scratch.PasteSpecial(-4122)   # xlPasteFormats

$scratch.Copy()
$dst.PasteSpecial(-4163)       # xlPasteValues
$scratch.Copy()
$dst.PasteSpecial(-4122)       # xlPasteFormats

$excel.CutCopyMode = 0

# Remove the scratch staging cells entirely (Clear(), not just
# ClearContents()+ClearFormats(), so the cells drop out of <sheetData>
# completely instead of lingering as empty shells).
$scratch.Clear()

# Remove the now-vacated M1 entirely (content + formatting) the same way.
$vacated = $ws.Range("M1")
$vacated.Clear()

# Update the active selection to match the post-edit state.
$ws.Range("C1:L1").Select() | Out-Null
